$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The practice table shrinks from 40 trial rows (A1:E41) to 20 trial rows
# (A1:E21): groups 1-4 now have 5 trials each instead of 10, and the colors
# used get reshuffled (yellow/green/purple/orange/brown/white -> magenta/
# darkviolet/yellow/saddlebrown/green). Delete the now-unused trailing rows
# first so row numbers/styles below line up with the final layout.
# ---------------------------------------------------------------------------
$ws.Rows("22:41").Delete()

# --- Group 1 (rows 2-6): left_color always red, right_color varies --------
$ws.Range("C2").Value = "magenta"
$ws.Range("C3").Value = "darkviolet"
$ws.Range("C4").Value = "yellow"
$ws.Range("C5").Value = "saddlebrown"
$ws.Range("C6").Value = "green"

# --- Group 2 (rows 7-11): both colors red, correct = "left,right" ---------
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "red"
$ws.Range("C7").Value = "red"
$ws.Range("D7").Value = "left,right"
$ws.Range("E7").Value = 2

$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "red"
$ws.Range("C8").Value = "red"
$ws.Range("D8").Value = "left,right"
$ws.Range("E8").Value = 2

$ws.Range("A9").Value = 13
$ws.Range("B9").Value = "red"
$ws.Range("C9").Value = "red"
$ws.Range("D9").Value = "left,right"
$ws.Range("E9").Value = 2

$ws.Range("A10").Value = 14
$ws.Range("B10").Value = "red"
$ws.Range("C10").Value = "red"
$ws.Range("D10").Value = "left,right"
$ws.Range("E10").Value = 2

$ws.Range("A11").Value = 15
$ws.Range("B11").Value = "red"
$ws.Range("C11").Value = "red"
$ws.Range("D11").Value = "left,right"
$ws.Range("E11").Value = 2

# --- Group 3 (rows 12-16): left_color varies, right_color always red ------
$ws.Range("A12").Value = 21
$ws.Range("B12").Value = "magenta"
$ws.Range("D12").Value = "right"
$ws.Range("E12").Value = 3

$ws.Range("A13").Value = 22
$ws.Range("B13").Value = "darkviolet"
$ws.Range("D13").Value = "right"
$ws.Range("E13").Value = 3

$ws.Range("A14").Value = 23
$ws.Range("B14").Value = "yellow"
$ws.Range("D14").Value = "right"
$ws.Range("E14").Value = 3

$ws.Range("A15").Value = 24
$ws.Range("B15").Value = "saddlebrown"
$ws.Range("D15").Value = "right"
$ws.Range("E15").Value = 3

$ws.Range("A16").Value = 25
$ws.Range("B16").Value = "green"
$ws.Range("D16").Value = "right"
$ws.Range("E16").Value = 3

# --- Group 4 (rows 17-21): both colors vary, different from each other ----
$ws.Range("A17").Value = 31
$ws.Range("B17").Value = "magenta"
$ws.Range("C17").Value = "darkviolet"
$ws.Range("D17").Value = "left,right,"
$ws.Range("E17").Value = 4

$ws.Range("A18").Value = 32
$ws.Range("B18").Value = "yellow"
$ws.Range("C18").Value = "green"
$ws.Range("D18").Value = "left,right,"
$ws.Range("E18").Value = 4

$ws.Range("A19").Value = 33
$ws.Range("B19").Value = "saddlebrown"
$ws.Range("C19").Value = "yellow"
$ws.Range("D19").Value = "left,right,"
$ws.Range("E19").Value = 4

$ws.Range("A20").Value = 34
$ws.Range("B20").Value = "yellow"
$ws.Range("C20").Value = "magenta"
$ws.Range("D20").Value = "left,right,"
$ws.Range("E20").Value = 4

$ws.Range("A21").Value = 35
$ws.Range("B21").Value = "green"
$ws.Range("C21").Value = "saddlebrown"
$ws.Range("D21").Value = "left,right,"
$ws.Range("E21").Value = 4

# ---------------------------------------------------------------------------
# Rows 7-11 and 17-21 used to belong to groups 1/4 (mismatched-color style:
# left_color border-set B, right_color border-set A) but now belong to
# groups 2/4 and 3 respectively - each group has its own highlight style for
# the left_color/right_color cells. Re-point those cells' formatting at an
# already-correctly-styled row from the matching group via a formats-only
# paste (keeps reusing the existing style indexes instead of cloning new
# ones).
# ---------------------------------------------------------------------------
$ws.Range("B12:C12").Copy()
$ws.Range("B7:C11").PasteSpecial(-4122)

$ws.Range("B3:C3").Copy()
$ws.Range("B17:C21").PasteSpecial(-4122)

$excel.CutCopyMode = 0
